$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1510096666666667
$ws.Range("H2").Value = 0.453029
$ws.Range("I2").Value = 0.01105950042918124
$ws.Range("J2").Value = 0.01105950042918124
$ws.Range("Q2").Value = 1.865886414702445
$ws.Range("R2").Value = 16.792977732322
$ws.Range("S2").Value = 0.01105950042918124
$ws.Range("T2").Value = 0.01105950042918124

# Row 3 updates
$ws.Range("I3").Value = 0.7495776481151314
$ws.Range("J3").Value = 0.7495776481151314
$ws.Range("S3").Value = 0.7495776481151314
$ws.Range("T3").Value = 0.7495776481151314

# Row 4 updates
$ws.Range("I4").Value = 0.2393628514556874
$ws.Range("J4").Value = 0.2393628514556874
$ws.Range("S4").Value = 0.2393628514556874
$ws.Range("T4").Value = 0.2393628514556874
